$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the dates in A2:A16 forward by exactly 365 days (one year later)
for ($r = 2; $r -le 16; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 365
}

# Update the selection to match the recorded view state (A2:A16, active cell A2)
$ws.Range("A2:A16").Select()
